# Applies the quantity (column F) corrections from the commit and
# recomputes all dependent values (column G line totals, per-company
# "Sub Total" rows, and the final "Sub Total"/"Grand Total" rows).
#
# Column layout on the sheet: A=Sl.No, B=Item Code, C=Item Name,
# D=Rate, E=Rate incl. tax, F=Qty, G=Value (= D * F).
# Each company block ends with a literal (non-formula) "Sub Total:" row
# whose B cell is the sum of the G column for that block; the very last
# two rows are a "Sub Total:"/"Grand Total:" pair whose B cells equal
# the sum of every per-company Sub Total B value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row/new-quantity pairs taken from the diff (column F).
$fChanges = @(
    @(20,342), @(21,99), @(26,57), @(28,59), @(36,123), @(41,243),
    @(46,63), @(52,104), @(55,141), @(58,97), @(71,3), @(227,74),
    @(306,33), @(386,5), @(392,148), @(397,99), @(413,108), @(416,85),
    @(432,137), @(433,188), @(450,15), @(453,52), @(454,99), @(466,33),
    @(511,324), @(527,80), @(532,38), @(534,162), @(544,66), @(614,86),
    @(617,53), @(620,413), @(621,29), @(626,36), @(659,65), @(673,45),
    @(674,1321), @(688,14), @(698,23)
)

# Update quantity and recompute the line value (G = D * F) for each
# changed row.
foreach ($chg in $fChanges) {
    $row = $chg[0]
    $newQty = $chg[1]
    $ws.Cells.Item($row, 6).Value = $newQty
    $rate = $ws.Cells.Item($row, 4).Value2
    $ws.Cells.Item($row, 7).Value = $rate * $newQty
}

# Per-company blocks that contain at least one changed row, described
# as (firstItemRow, lastItemRow, subTotalRow).
$affectedGroups = @(
    @(18,33,34), @(36,65,66), @(71,82,83), @(220,228,229),
    @(297,327,328), @(380,388,389), @(391,416,417), @(429,437,438),
    @(450,457,458), @(465,475,476), @(510,524,525), @(527,534,535),
    @(543,555,556), @(609,627,628), @(659,667,668), @(673,679,680),
    @(686,690,691), @(693,712,713)
)

foreach ($grp in $affectedGroups) {
    $first = $grp[0]
    $last = $grp[1]
    $subRow = $grp[2]
    $sum = 0
    for ($r = $first; $r -le $last; $r++) {
        $sum = $sum + $ws.Cells.Item($r, 7).Value2
    }
    $ws.Cells.Item($subRow, 2).Value = $sum
}

# Every per-company "Sub Total:" row on the sheet (excluding the final
# aggregate pair) so the grand-total rows can be rebuilt from scratch.
$allSubTotalRows = @(
    9,12,16,34,66,69,83,87,90,97,123,128,133,137,147,155,164,170,193,
    197,201,204,208,218,229,240,248,295,328,332,336,340,356,363,372,
    375,378,389,417,427,438,441,445,448,458,463,476,482,493,508,525,
    535,538,541,556,561,573,582,585,588,595,603,607,628,635,657,668,
    671,680,684,691,713,717
)

$grandTotal = 0
foreach ($r in $allSubTotalRows) {
    $grandTotal = $grandTotal + $ws.Cells.Item($r, 2).Value2
}

# Row 718 is labelled "Sub Total:" (sum of all sub totals) and row 719
# is the "Grand Total:" row, which mirrors the same value.
$ws.Cells.Item(718, 2).Value = $grandTotal
$ws.Cells.Item(719, 2).Value = $grandTotal
